$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Update the "All" column counts/percentages: the analysis now
#    only includes the insured (Private Insurance) and uninsured
#    (Self-pay) subgroups, so the combined "All" totals shrink.
# ------------------------------------------------------------------
$d.Content.Find.Execute('All (n=274989)', $true, $false, $false, $false, $false, $true, 1, $false, 'All (n=197870)', 2) | Out-Null
$d.Content.Find.Execute('40.34±13.35', $true, $false, $false, $false, $false, $true, 1, $false, '40.59±13.27', 2) | Out-Null
$d.Content.Find.Execute('9480 (3.45)', $true, $false, $false, $false, $false, $true, 1, $false, '6830 (3.45)', 2) | Out-Null
$d.Content.Find.Execute('23574 (8.57)', $true, $false, $false, $false, $false, $true, 1, $false, '14143 (7.15)', 2) | Out-Null
$d.Content.Find.Execute('55990 (20.36)', $true, $false, $false, $false, $false, $true, 1, $false, '33054 (16.70)', 2) | Out-Null
$d.Content.Find.Execute('1575 (0.57)', $true, $false, $false, $false, $false, $true, 1, $false, '891 (0.45)', 2) | Out-Null
$d.Content.Find.Execute('11339 (4.12)', $true, $false, $false, $false, $false, $true, 1, $false, '7614 (3.85)', 2) | Out-Null
$d.Content.Find.Execute('15995 (5.82)', $true, $false, $false, $false, $false, $true, 1, $false, '12130 (6.13)', 2) | Out-Null
$d.Content.Find.Execute('157036 (57.11)', $true, $false, $false, $false, $false, $true, 1, $false, '123208 (62.27)', 2) | Out-Null
$d.Content.Find.Execute('131408 (47.79)', $true, $false, $false, $false, $false, $true, 1, $false, '92056 (46.52)', 2) | Out-Null
$d.Content.Find.Execute('142467 (51.81)', $true, $false, $false, $false, $false, $true, 1, $false, '104972 (53.05)', 2) | Out-Null
$d.Content.Find.Execute('1114 (0.41)', $true, $false, $false, $false, $false, $true, 1, $false, '842 (0.43)', 2) | Out-Null
$d.Content.Find.Execute('71903 (26.15)', $true, $false, $false, $false, $false, $true, 1, $false, '44420 (22.45)', 2) | Out-Null
$d.Content.Find.Execute('67128 (24.41)', $true, $false, $false, $false, $false, $true, 1, $false, '46239 (23.37)', 2) | Out-Null
$d.Content.Find.Execute('69846 (25.40)', $true, $false, $false, $false, $false, $true, 1, $false, '52340 (26.45)', 2) | Out-Null
$d.Content.Find.Execute('66112 (24.04)', $true, $false, $false, $false, $false, $true, 1, $false, '54871 (27.73)', 2) | Out-Null
$d.Content.Find.Execute('27877 (10.14)', $true, $false, $false, $false, $false, $true, 1, $false, '19918 (10.07)', 2) | Out-Null
$d.Content.Find.Execute('103415 (37.61)', $true, $false, $false, $false, $false, $true, 1, $false, '76850 (38.84)', 2) | Out-Null
$d.Content.Find.Execute('143697 (52.26)', $true, $false, $false, $false, $false, $true, 1, $false, '101102 (51.10)', 2) | Out-Null
$d.Content.Find.Execute('48168 (17.52)', $true, $false, $false, $false, $false, $true, 1, $false, '36571 (18.48)', 2) | Out-Null
$d.Content.Find.Execute('56761 (20.64)', $true, $false, $false, $false, $false, $true, 1, $false, '39610 (20.02)', 2) | Out-Null
$d.Content.Find.Execute('97440 (35.43)', $true, $false, $false, $false, $false, $true, 1, $false, '74508 (37.66)', 2) | Out-Null
$d.Content.Find.Execute('72620 (26.41)', $true, $false, $false, $false, $false, $true, 1, $false, '47181 (23.84)', 2) | Out-Null
$d.Content.Find.Execute('229819 (83.57)', $true, $false, $false, $false, $false, $true, 1, $false, '169619 (85.72)', 2) | Out-Null
$d.Content.Find.Execute('40947 (14.89)', $true, $false, $false, $false, $false, $true, 1, $false, '25434 (12.85)', 2) | Out-Null
$d.Content.Find.Execute('4223 (1.54)', $true, $false, $false, $false, $false, $true, 1, $false, '2817 (1.42)', 2) | Out-Null
$d.Content.Find.Execute('150145 (54.60)', $true, $false, $false, $false, $false, $true, 1, $false, '110524 (55.86)', 2) | Out-Null
$d.Content.Find.Execute('107709 (39.17)', $true, $false, $false, $false, $false, $true, 1, $false, '77406 (39.12)', 2) | Out-Null
$d.Content.Find.Execute('14885 (5.41)', $true, $false, $false, $false, $false, $true, 1, $false, '8690 (4.39)', 2) | Out-Null
$d.Content.Find.Execute('2250 (0.82)', $true, $false, $false, $false, $false, $true, 1, $false, '1250 (0.63)', 2) | Out-Null
$d.Content.Find.Execute('251390 (91.42)', $true, $false, $false, $false, $false, $true, 1, $false, '184177 (93.08)', 2) | Out-Null
$d.Content.Find.Execute('17539 (6.38)', $true, $false, $false, $false, $false, $true, 1, $false, '10446 (5.28)', 2) | Out-Null
$d.Content.Find.Execute('4865 (1.77)', $true, $false, $false, $false, $false, $true, 1, $false, '2585 (1.31)', 2) | Out-Null
$d.Content.Find.Execute('1195 (0.43)', $true, $false, $false, $false, $false, $true, 1, $false, '662 (0.33)', 2) | Out-Null
$d.Content.Find.Execute('0.33±1.15', $true, $false, $false, $false, $false, $true, 1, $false, '0.25±1.00', 2) | Out-Null
$d.Content.Find.Execute('213 (0.08)', $true, $false, $false, $false, $false, $true, 1, $false, '88 (0.04)', 2) | Out-Null
$d.Content.Find.Execute('68404 (24.88)', $true, $false, $false, $false, $false, $true, 1, $false, '46885 (23.69)', 2) | Out-Null
$d.Content.Find.Execute('25110 (9.13)', $true, $false, $false, $false, $false, $true, 1, $false, '16676 (8.43)', 2) | Out-Null
$d.Content.Find.Execute('1187 (0.43)', $true, $false, $false, $false, $false, $true, 1, $false, '840 (0.42)', 2) | Out-Null

# ------------------------------------------------------------------
# 2) Insert a new "CCI > 0" row right after the "CCI Score" row.
# ------------------------------------------------------------------
$t = $d.Tables.Item(1)
$cciScoreRow = 0
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    if ($t.Cell($i, 1).Range.Text -like "CCI Score*") {
        $cciScoreRow = $i
        break
    }
}

$newRow = $t.Rows.Add($t.Rows.Item($cciScoreRow + 1))
$t.Cell($cciScoreRow + 1, 1).Range.Text = "CCI > 0"
$t.Cell($cciScoreRow + 1, 2).Range.Text = "27857 (14.08)"
$t.Cell($cciScoreRow + 1, 3).Range.Text = "23480 (14.82)"
$t.Cell($cciScoreRow + 1, 4).Range.Text = "4377 (11.09)"
